$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 4), matching the date style already used in A2:A3
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A4").Value = 44313
$ws.Range("B4").Value = 6.5
$ws.Range("C4").Value = "Сделалсохранение в локал стор, добавил валидацию на странице логина."

# Recalculate so the SUM(B:B) formula in F2 reflects the new row
$wb.Application.Calculate()
